# Weekly update: insert a new observation row for "Jengibre" at row 37,
# pushing the existing rows 37:70 down to 38:71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 37 (shifts 37:70 -> 38:71,
# and the sheet's used-range / dimension grows to A1:R71 automatically).
$ws.Rows("37").Insert()

# Populate the newly inserted row 37 with the new weekly record.
$ws.Range("A37").Value = 6
$ws.Range("B37").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C37").Value = "Metropolitana"
$ws.Range("D37").Value = 44651
$ws.Range("E37").Value = 13
$ws.Range("F37").Value = 100114007
$ws.Range("G37").Value = "Jengibre"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 180
$ws.Range("K37").Value = 12000
$ws.Range("L37").Value = 14000
$ws.Range("M37").Value = 12889
$ws.Range("N37").Value = '$/caja 13 kilos'
$ws.Range("O37").Value = "Perú"
$ws.Range("P37").Value = 991
$ws.Range("Q37").Value = 13
$ws.Range("R37").Value = "Hortaliza"
